# Applies the "Nuevos modelos de ML" edit: appends new microciclo rows
# (146-153) to the Microciclos sheet, including a new "Mazatlan" PARTIDO
# entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New rows to append after the existing data (row 145).
# Columns: A=Fecha(date serial), B=Microciclo_Num, C=Tipo_Microciclo,
#          D=Tipo_Dia, E=Fase, F=Intensidad, G=Partido
$rows = @(
    @{ Row = 146; Fecha = 45942; Num = 21; Tipo = "Competencia"; Dia = "Competencia"; Fase = "DESCANSO" },
    @{ Row = 147; Fecha = 45943; Num = 22; Tipo = "Competencia"; Dia = "Competencia"; Fase = "DESCANSO" },
    @{ Row = 148; Fecha = 45944; Num = 22; Tipo = "Competencia"; Dia = "Competencia"; Fase = "ENTRENO"; Intensidad = 1 },
    @{ Row = 149; Fecha = 45945; Num = 22; Tipo = "Competencia"; Dia = "Competencia"; Fase = "ENTRENO"; Intensidad = 2 },
    @{ Row = 150; Fecha = 45946; Num = 22; Tipo = "Competencia"; Dia = "Competencia"; Fase = "ENTRENO"; Intensidad = -2 },
    @{ Row = 151; Fecha = 45947; Num = 22; Tipo = "Competencia"; Dia = "Competencia"; Fase = "ENTRENO"; Intensidad = -1 },
    @{ Row = 152; Fecha = 45948; Num = 22; Tipo = "Competencia"; Dia = "Competencia"; Fase = "PARTIDO"; Partido = "Mazatlan" },
    @{ Row = 153; Fecha = 45949; Num = 22; Tipo = "Competencia"; Dia = "Competencia"; Fase = "ENTRENO"; Intensidad = 1 }
)

foreach ($r in $rows) {
    $rowNum = $r.Row

    $cellA = $ws.Cells.Item($rowNum, 1)
    $cellA.Value = $r.Fecha
    # Reuse the existing builtin date style (numFmtId 14, same as the rest
    # of column A) instead of a format string that would register a brand
    # new custom numFmt entry in styles.xml.
    $cellA.NumberFormat = "m/d/yy"

    $ws.Cells.Item($rowNum, 2).Value = $r.Num
    $ws.Cells.Item($rowNum, 3).Value = $r.Tipo
    $ws.Cells.Item($rowNum, 4).Value = $r.Dia
    $ws.Cells.Item($rowNum, 5).Value = $r.Fase

    if ($r.ContainsKey("Intensidad")) {
        $ws.Cells.Item($rowNum, 6).Value = $r.Intensidad
    }
    if ($r.ContainsKey("Partido")) {
        $ws.Cells.Item($rowNum, 7).Value = $r.Partido
    }
}

# Update the view to match the final saved state (scrolled down, new
# active selection below the appended data).
$ws.Application.ActiveWindow.ScrollRow = 142
$ws.Range("A154").Select()
